$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A (probe ids) for rows 2-6 first
$ws.Range("A2").Value = "cg21922223"
$ws.Range("A3").Value = "cg23719692"
$ws.Range("A4").Value = "cg10678427"
$ws.Range("A5").Value = "cg11189107"
$ws.Range("A6").Value = "cg01877778"

# Then update column B (gene / region) for rows 2-6
$ws.Range("B2").Value = "non-genic"
$ws.Range("B3").Value = "AGAP1"
$ws.Range("B4").Value = "non-genic"
$ws.Range("B5").Value = "CTSG"
$ws.Range("B6").Value = "PTPRN2"

# Update the selected cell to match the new active selection
$ws.Range("B6").Select()
